# Apply updated cryptocurrency price/volume data to the worksheet
# (cryptos list refresh, per GitHub Actions scheduled update)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.901.94'
$ws.Range('E2').Value = '  +2.68%  '
$ws.Range('D3').Value = '2.113.69'
$ws.Range('E3').Value = '  +10.13%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.005'
$ws.Range('E4').Value = '  +0.23%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '335.32'
$ws.Range('E5').Value = '  +5.03%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.003'
$ws.Range('E6').Value = '  +0.15%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5309'
$ws.Range('E7').Value = '  +4.61%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4369'
$ws.Range('E8').Value = '  +8.57%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.09034'
$ws.Range('E9').Value = '  +8.61%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '46.26'
$ws.Range('E10').Value = '  +9.81%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.177'
$ws.Range('E11').Value = '  +5.67%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '25.12'
$ws.Range('E12').Value = '  +4.22%  '
$ws.Range('D13').Value = '2.117.01'
$ws.Range('E13').Value = '  +10.39%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.777'
$ws.Range('E14').Value = '  +5.63%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.805'
$ws.Range('E15').Value = '  +7.74%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '97.59'
$ws.Range('E16').Value = '  +5.53%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.004'
$ws.Range('E17').Value = '  +0.33%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001134'
$ws.Range('E18').Value = '  +3.60%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06669'
$ws.Range('E19').Value = '  +2.71%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '19.14'
$ws.Range('E20').Value = '  +3.80%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.002'
$ws.Range('E21').Value = '  +0.05%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.368'
$ws.Range('E22').Value = '  +7.11%  '
$ws.Range('D23').Value = '30.965.93'
$ws.Range('E23').Value = '  +2.89%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '12.16'
$ws.Range('E24').Value = '  +7.27%  '
$ws.Range('D25').Value = '2.365.87'
$ws.Range('E25').Value = '  +10.71%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.274'
$ws.Range('E26').Value = '  +3.60%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '22.79'
$ws.Range('E27').Value = '  +4.49%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.566'
$ws.Range('E28').Value = '  +13.01%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '163.17'
$ws.Range('E29').Value = '  +0.42%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '133.87'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.170'
$ws.Range('E31').Value = '  +3.28%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.1073'
$ws.Range('E32').Value = '  +2.60%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.229'
$ws.Range('E33').Value = '  +4.88%  '
$ws.Range('E34').Value = '  +6.07%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.521'
$ws.Range('E35').Value = '  +23.03%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02615'
$ws.Range('E36').Value = '  +6.78%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '12.96'
$ws.Range('E37').Value = '  +12.74%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.545'
$ws.Range('E38').Value = '  +4.51%  '
# Rows 39/40: Hedera and FraxShare swapped ranking order
$ws.Range('B39').Value = 'FraxShare'
$ws.Range('C39').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '9.525'
$ws.Range('E39').Value = '  +10.48%  '
$ws.Range('B40').Value = 'Hedera'
$ws.Range('C40').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.06724'
$ws.Range('E40').Value = '  +4.35%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.2276'
$ws.Range('E41').Value = '  +6.08%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.6866'
$ws.Range('E42').Value = '  +6.25%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.252'
$ws.Range('E43').Value = '  +3.10%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.6451'
$ws.Range('E44').Value = '  +6.77%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '14.13'
$ws.Range('E45').Value = '  +6.14%  '
$ws.Range('E46').Value = '  +0.15%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.248'
$ws.Range('E47').Value = '  +3.81%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.682'
$ws.Range('E48').Value = '  +1.71%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.277'
$ws.Range('E49').Value = '  +5.83%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '82.87'
$ws.Range('E50').Value = '  +6.40%  '
# Row 51: Cronos replaced by WEMIXTOKEN
$ws.Range('B51').Value = 'WEMIXTOKEN'
$ws.Range('C51').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.169'
$ws.Range('E51').Value = '  +2.85%  '
